$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$data = @(
    @("TestCase_B15", "To verify that 10 MORE button is working correctly", "Y", "FAIL"),
    @("TestCase_B16", "To verify that 10 MORE button is not present in search results page if the total search results is less than or equal to 10", "Y", "PASS"),
    @("TestCase_B17", "To verify that 10 MORE button is present in search results page if total search results is more than 10", "Y", "PASS"),
    @("TestCase_B18", "To verify that Times cited,Views,Comments fields are getting displayed for each document in search results page", "Y", "PASS"),
    @("TestCase_B19", "To verify that user is able to expand and collapse SORT BY drop down", "Y", "PASS"),
    @("TestCase_B20", "To verify that user is able to sort the documents by TIMES CITED field", "Y", "PASS"),
    @("TestCase_B21", "To verify that RESET button in the left navigation pane in search results page is working correctly", "Y", "PASS"),
    @("TestCase_B22", "To verify that MORE and LESS links in the left navigation pane are working correctly", "Y", "PASS"),
    @("TestCase_B23", "To verify that user is able to collapse and expand the filters in the left navigation pane in search results page", "Y", "FAIL")
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]

    $ws.Range("A3").Copy()
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 4)).PasteSpecial(-4122)
}

$ws.Range("C18").Select()
